$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6678836941719055
$ws.Range("B1").Value = 0.9364299774169922
$ws.Range("C1").Value = 1.221186518669128
$ws.Range("D1").Value = 3.882378339767456
$ws.Range("E1").Value = 2.425849199295044
